{"js": "// Update the Data Request Letter:\n// Insert \"let us know and \" right before \"we will honor that request.\"\n// so the sentence reads \"...publicly available, then let us know and\n// we will honor that request.\"\n\nconst body = context.document.body;\n\n// Anchor on a long, unique phrase so the insertion point is unambiguous.\nconst searchResults = body.search(\"we will honor that request.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the target phrase \"we will honor that request.\"');\n}\n\nconst target = searchResults.items[0];\ntarget.insertText(\"let us know and \", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Update the Data Request Letter:\n# Insert \"let us know and \" right before \"we will honor that request.\"\n# so the sentence reads \"...publicly available, then let us know and\n# we will honor that request.\"\n\n$d = $word.ActiveDocument\n\n# Anchor on a long, unique phrase so the insertion point is unambiguous.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"we will honor that request.\"\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop\n$find.Execute() | Out-Null\n\nif (-not $find.Found) {\n    throw 'Could not find the target phrase \"we will honor that request.\"'\n}\n\n$target = $find.Parent\n# Collapse to the start of the found range so the insertion lands\n# immediately before \"we will honor that request.\"\n$insertionRange = $d.Range($target.Start, $target.Start)\n$insertionRange.InsertBefore(\"let us know and \")\n"}
